# Tskaltubo Municipality - "Area" sheet: revert the 3-census-year table
# (1989 / 2002 / 2014, columns B:D) back to the single-year 2014 layout
# (column B only), restoring the simpler title/subtitle arrangement.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Drop the historical census columns (1989, 2002). The remaining "2014"
# column (old D) slides into column B, carrying its own number format /
# borders with it.
$ws.Columns("B:C").Delete()

# Drop the blank spacer row that used to sit between the "Area" label and
# the "(according to the population census data)" subtitle.
$ws.Rows("3").Delete()

# "Area" used to live right under the title (row 2); the sheet no longer
# shows it there - it now labels the data row instead (see below), so
# clear it out and leave row 2 blank.
$ws.Range("A2").Clear()

# What used to be "(according to the population census data)" is replaced
# by the simpler "(sq. km)" unit caption.
$ws.Range("A3").Value = "(sq. km)"

# The row that used to carry the "(sq. km)" row label now carries "Area".
$ws.Range("A5").Value = "Area"

# Restore the original (taller) row heights used by the single-year layout.
$ws.Rows("1:7").RowHeight = 20.1

$wb.Save()
